# Fruta / hortaliza, semanal
# Insert a new weekly record at the top of the data (row 17), pushing the
# existing rows 17-26 down to 18-27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17 (shifts rows 17:26 down to 18:27,
# inheriting formatting - e.g. the date style on column D - from the row
# that used to be there).
$ws.Rows.Item(17).Insert()

# Populate the newly inserted row 17 with the new weekly observation.
$ws.Cells.Item(17, 1).Value() = 1
$ws.Cells.Item(17, 2).Value() = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(17, 3).Value() = "Arica y Parinacota"
$ws.Cells.Item(17, 4).Value() = 45096
$ws.Cells.Item(17, 5).Value() = 15
$ws.Cells.Item(17, 6).Value() = 100114007
$ws.Cells.Item(17, 7).Value() = "Jengibre"
$ws.Cells.Item(17, 8).Value() = "Sin especificar"
$ws.Cells.Item(17, 9).Value() = "Primera"
$ws.Cells.Item(17, 10).Value() = 750
$ws.Cells.Item(17, 11).Value() = 14000
$ws.Cells.Item(17, 12).Value() = 15000
$ws.Cells.Item(17, 13).Value() = 14600
$ws.Cells.Item(17, 14).Value() = "$/caja 13 kilos"
$ws.Cells.Item(17, 15).Value() = "Perú"
$ws.Cells.Item(17, 16).Value() = 1123
$ws.Cells.Item(17, 17).Value() = 13
$ws.Cells.Item(17, 18).Value() = "Hortaliza"
